$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subsystem Progress Tracking")
$ws.Columns.Item(8).Delete()
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(4).Delete()
Write-Host $ws.UsedRange.Address()
